{"js": "// Update the date line and every three-digit-by-one-digit multiplication\n// answer in the practice table to the new day's problems/answers.\nconst replacements = [\n  [\"2024-02-13 Tuesday\", \"2024-02-14 Wednesday\"],\n  [\"625\u00d75=3125\", \"611\u00d78=4888\"],\n  [\"815\u00d78=6520\", \"539\u00d72=1078\"],\n  [\"223\u00d73=669\", \"212\u00d78=1696\"],\n  [\"590\u00d73=1770\", \"711\u00d76=4266\"],\n  [\"508\u00d72=1016\", \"818\u00d78=6544\"],\n  [\"393\u00d79=3537\", \"470\u00d78=3760\"],\n  [\"756\u00d75=3780\", \"947\u00d75=4735\"],\n  [\"566\u00d74=2264\", \"715\u00d73=2145\"],\n  [\"404\u00d77=2828\", \"725\u00d76=4350\"],\n  [\"456\u00d75=2280\", \"692\u00d74=2768\"],\n  [\"983\u00d79=8847\", \"305\u00d72=610\"],\n  [\"553\u00d78=4424\", \"443\u00d79=3987\"],\n  [\"828\u00d75=4140\", \"420\u00d75=2100\"],\n  [\"791\u00d75=3955\", \"885\u00d79=7965\"],\n  [\"942\u00d76=5652\", \"984\u00d72=1968\"],\n  [\"273\u00d72=546\", \"529\u00d77=3703\"],\n  [\"500\u00d72=1000\", \"340\u00d72=680\"],\n  [\"647\u00d79=5823\", \"309\u00d74=1236\"],\n  [\"263\u00d74=1052\", \"107\u00d72=214\"],\n  [\"803\u00d78=6424\", \"436\u00d77=3052\"],\n  [\"865\u00d77=6055\", \"768\u00d77=5376\"],\n  [\"679\u00d77=4753\", \"628\u00d74=2512\"],\n  [\"874\u00d75=4370\", \"852\u00d76=5112\"],\n  [\"638\u00d74=2552\", \"132\u00d79=1188\"],\n  [\"449\u00d79=4041\", \"185\u00d75=925\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the date line and every three-digit-by-one-digit multiplication\n# answer in the practice table to the new day's problems/answers.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2024-02-13 Tuesday\", \"2024-02-14 Wednesday\"),\n    @(\"625\u00d75=3125\", \"611\u00d78=4888\"),\n    @(\"815\u00d78=6520\", \"539\u00d72=1078\"),\n    @(\"223\u00d73=669\", \"212\u00d78=1696\"),\n    @(\"590\u00d73=1770\", \"711\u00d76=4266\"),\n    @(\"508\u00d72=1016\", \"818\u00d78=6544\"),\n    @(\"393\u00d79=3537\", \"470\u00d78=3760\"),\n    @(\"756\u00d75=3780\", \"947\u00d75=4735\"),\n    @(\"566\u00d74=2264\", \"715\u00d73=2145\"),\n    @(\"404\u00d77=2828\", \"725\u00d76=4350\"),\n    @(\"456\u00d75=2280\", \"692\u00d74=2768\"),\n    @(\"983\u00d79=8847\", \"305\u00d72=610\"),\n    @(\"553\u00d78=4424\", \"443\u00d79=3987\"),\n    @(\"828\u00d75=4140\", \"420\u00d75=2100\"),\n    @(\"791\u00d75=3955\", \"885\u00d79=7965\"),\n    @(\"942\u00d76=5652\", \"984\u00d72=1968\"),\n    @(\"273\u00d72=546\", \"529\u00d77=3703\"),\n    @(\"500\u00d72=1000\", \"340\u00d72=680\"),\n    @(\"647\u00d79=5823\", \"309\u00d74=1236\"),\n    @(\"263\u00d74=1052\", \"107\u00d72=214\"),\n    @(\"803\u00d78=6424\", \"436\u00d77=3052\"),\n    @(\"865\u00d77=6055\", \"768\u00d77=5376\"),\n    @(\"679\u00d77=4753\", \"628\u00d74=2512\"),\n    @(\"874\u00d75=4370\", \"852\u00d76=5112\"),\n    @(\"638\u00d74=2552\", \"132\u00d79=1188\"),\n    @(\"449\u00d79=4041\", \"185\u00d75=925\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
